$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.826.72"
$ws.Range("E2").Value = "  +0.76%  "

$ws.Range("D3").Value = "3.532.88"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.14"
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.01"
$ws.Range("E6").Value = "  +5.86%  "

$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.203"
$ws.Range("E9").Value = "  -3.22%  "

$ws.Range("E10").Value = "  -0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.68"
$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("D14").Value = "4.092.23"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "603.58"
$ws.Range("E15").Value = "  -1.28%  "

$ws.Range("D16").Value = "69.990.42"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.12"
$ws.Range("E17").Value = "  +1.18%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.75"

$ws.Range("D19").Value = "3.538.37"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("E21").Value = "  +0.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.24"
$ws.Range("E22").Value = "  +4.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +5.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "102.29"
$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("E25").Value = "  -0.45%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.14"
$ws.Range("E26").Value = "  +4.74%  "

$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.49"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.33"
$ws.Range("E30").Value = "  +17.22%  "

$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.56"
$ws.Range("E32").Value = "  +1.57%  "

$ws.Range("E33").Value = "  -0.36%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.23"
$ws.Range("E34").Value = "  -0.14%  "

$ws.Range("E35").Value = "  +12.63%  "

$ws.Range("D36").Value = "3.741.92"
$ws.Range("E36").Value = "  +5.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  -2.55%  "

$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.63"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "485.59"
$ws.Range("E42").Value = "  -7.74%  "

$ws.Range("E43").Value = "  -4.86%  "

$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("E45").Value = "  -2.70%  "

$ws.Range("E46").Value = "  -1.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.31"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("E49").Value = "  -3.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000253"
$ws.Range("E50").Value = "  +6.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.72"
$ws.Range("E51").Value = "  -0.70%  "
